# Applies the weekly Fruta/Hortalizas price-sheet update for
# "Terminal Hortofrutícola Agro Chillán - Membrillo":
#  - revises the price/volume/quality/date/origin data recorded in rows 20-34
#  - appends three brand-new observation rows (35-37)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 20-34 with revised weekly price data ---
# Row 20
$ws.Range("D20").Value = 45090
$ws.Range("L20").Value = 'Especial'
$ws.Range("M20").Value = 80
$ws.Range("N20").Value = 11000
$ws.Range("O20").Value = 11000
$ws.Range("P20").Value = 11000
$ws.Range("Q20").Value = '$/caja 18 kilos empedrada'
$ws.Range("R20").Value = 'Región del Maule'
$ws.Range("S20").Value = 611

# Row 21
$ws.Range("D21").Value = 45090
$ws.Range("L21").Value = 'Primera'
$ws.Range("N21").Value = 10000
$ws.Range("O21").Value = 10000
$ws.Range("P21").Value = 10000
$ws.Range("R21").Value = 'Región del Maule'
$ws.Range("S21").Value = 556

# Row 22
$ws.Range("D22").Value = 45090
$ws.Range("L22").Value = 'Segunda'
$ws.Range("N22").Value = 8000
$ws.Range("O22").Value = 8000
$ws.Range("P22").Value = 8000
$ws.Range("R22").Value = 'Región del Maule'
$ws.Range("S22").Value = 444

# Row 23
$ws.Range("D23").Value = 45020
$ws.Range("L23").Value = 'Primera'
$ws.Range("N23").Value = 12000
$ws.Range("O23").Value = 12000
$ws.Range("P23").Value = 12000
$ws.Range("Q23").Value = '$/caja 18 kilos granel'
$ws.Range("R23").Value = 'Región de O''Higgins'
$ws.Range("S23").Value = 667

# Row 24
$ws.Range("D24").Value = 45040
$ws.Range("L24").Value = 'Especial'
$ws.Range("N24").Value = 13000
$ws.Range("O24").Value = 13000
$ws.Range("P24").Value = 13000
$ws.Range("R24").Value = 'Región de O''Higgins'
$ws.Range("S24").Value = 722

# Row 25
$ws.Range("D25").Value = 45040
$ws.Range("L25").Value = 'Primera'
$ws.Range("M25").Value = 40
$ws.Range("N25").Value = 12000
$ws.Range("O25").Value = 12000
$ws.Range("P25").Value = 12000
$ws.Range("R25").Value = 'Región de O''Higgins'
$ws.Range("S25").Value = 667

# Row 26
$ws.Range("D26").Value = 45089
$ws.Range("L26").Value = 'Especial'
$ws.Range("N26").Value = 11000
$ws.Range("O26").Value = 11000
$ws.Range("P26").Value = 11000
$ws.Range("R26").Value = 'Región del Maule'
$ws.Range("S26").Value = 611

# Row 27
$ws.Range("D27").Value = 45089
$ws.Range("L27").Value = 'Primera'
$ws.Range("N27").Value = 9000
$ws.Range("O27").Value = 9000
$ws.Range("P27").Value = 9000
$ws.Range("R27").Value = 'Región del Maule'
$ws.Range("S27").Value = 500

# Row 28
$ws.Range("D28").Value = 45089
$ws.Range("L28").Value = 'Segunda'
$ws.Range("M28").Value = 30
$ws.Range("N28").Value = 7000
$ws.Range("O28").Value = 7000
$ws.Range("P28").Value = 7000
$ws.Range("R28").Value = 'Región del Maule'
$ws.Range("S28").Value = 389

# Row 29
$ws.Range("D29").Value = 45070
$ws.Range("M29").Value = 60
$ws.Range("R29").Value = 'Región de O''Higgins'

# Row 30
$ws.Range("D30").Value = 45062
$ws.Range("M30").Value = 50

# Row 31
$ws.Range("D31").Value = 45062
$ws.Range("M31").Value = 50

# Row 32
$ws.Range("D32").Value = 45085
$ws.Range("M32").Value = 50
$ws.Range("N32").Value = 10000
$ws.Range("O32").Value = 10000
$ws.Range("P32").Value = 10000
$ws.Range("Q32").Value = '$/caja 18 kilos empedrada'
$ws.Range("R32").Value = 'Región del Maule'
$ws.Range("S32").Value = 556
$ws.Range("T32").Value = 18

# Row 33
$ws.Range("D33").Value = 45033
$ws.Range("L33").Value = 'Especial'
$ws.Range("M33").Value = 60
$ws.Range("N33").Value = 13000
$ws.Range("O33").Value = 13000
$ws.Range("P33").Value = 13000
$ws.Range("Q33").Value = '$/caja 18 kilos empedrada'
$ws.Range("S33").Value = 722
$ws.Range("T33").Value = 18

# Row 34
$ws.Range("D34").Value = 45033
$ws.Range("M34").Value = 80
$ws.Range("Q34").Value = '$/caja 18 kilos empedrada'

# --- Append brand-new rows 35-37 ---
# Row 35
$ws.Range("A35").Value = 7
$ws.Range("B35").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C35").Value = 'Ñuble'
$ws.Range("D35").Value = 45076
$ws.Range("D35").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E35").Value = 16
$ws.Range("F35").Value = 'Fruta'
$ws.Range("G35").Value = 100104
$ws.Range("H35").Value = 'Frutos de pepita'
$ws.Range("I35").Value = 100104003
$ws.Range("J35").Value = 'Membrillo'
$ws.Range("K35").Value = 'Champion'
$ws.Range("L35").Value = 'Primera'
$ws.Range("M35").Value = 30
$ws.Range("N35").Value = 12000
$ws.Range("O35").Value = 12000
$ws.Range("P35").Value = 12000
$ws.Range("Q35").Value = '$/caja 15 kilos granel'
$ws.Range("R35").Value = 'Región de O''Higgins'
$ws.Range("S35").Value = 800
$ws.Range("T35").Value = 15

# Row 36
$ws.Range("A36").Value = 7
$ws.Range("B36").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C36").Value = 'Ñuble'
$ws.Range("D36").Value = 45076
$ws.Range("D36").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E36").Value = 16
$ws.Range("F36").Value = 'Fruta'
$ws.Range("G36").Value = 100104
$ws.Range("H36").Value = 'Frutos de pepita'
$ws.Range("I36").Value = 100104003
$ws.Range("J36").Value = 'Membrillo'
$ws.Range("K36").Value = 'Champion'
$ws.Range("L36").Value = 'Segunda'
$ws.Range("M36").Value = 30
$ws.Range("N36").Value = 10000
$ws.Range("O36").Value = 10000
$ws.Range("P36").Value = 10000
$ws.Range("Q36").Value = '$/caja 15 kilos granel'
$ws.Range("R36").Value = 'Región de O''Higgins'
$ws.Range("S36").Value = 667
$ws.Range("T36").Value = 15

# Row 37
$ws.Range("A37").Value = 7
$ws.Range("B37").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C37").Value = 'Ñuble'
$ws.Range("D37").Value = 45021
$ws.Range("D37").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E37").Value = 16
$ws.Range("F37").Value = 'Fruta'
$ws.Range("G37").Value = 100104
$ws.Range("H37").Value = 'Frutos de pepita'
$ws.Range("I37").Value = 100104003
$ws.Range("J37").Value = 'Membrillo'
$ws.Range("K37").Value = 'Champion'
$ws.Range("L37").Value = 'Primera'
$ws.Range("M37").Value = 50
$ws.Range("N37").Value = 12000
$ws.Range("O37").Value = 12000
$ws.Range("P37").Value = 12000
$ws.Range("Q37").Value = '$/caja 18 kilos granel'
$ws.Range("R37").Value = 'Región de O''Higgins'
$ws.Range("S37").Value = 667
$ws.Range("T37").Value = 18
